$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.382.25'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = '3.674.83'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = "'643.13"
$ws.Range("E5").Value = '  -5.32%  '
$ws.Range("D6").Value = "'159.88"
$ws.Range("E6").Value = '  +0.51%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +0.58%  '
$ws.Range("E9").Value = '  -0.35%  '
$ws.Range("E10").Value = '  -0.60%  '
$ws.Range("D11").Value = "'0.449"
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("D13").Value = '4.294.06'
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("D14").Value = "'32.71"
$ws.Range("E14").Value = '  +0.89%  '
$ws.Range("D15").Value = '3.647.28'
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").Value = '69.357.50'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = "'16.03"
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("D19").Value = "'6.49"
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").Value = "'466.20"
$ws.Range("E20").Value = '  -0.23%  '
$ws.Range("D21").Value = "'9.92"
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").Value = "'0.648"
$ws.Range("E22").Value = '  -0.97%  '
$ws.Range("E23").Value = '  -0.66%  '
$ws.Range("D24").Value = '3.820.06'
$ws.Range("E24").Value = '  -0.38%  '
$ws.Range("E26").Value = '  +3.25%  '
$ws.Range("D27").Value = "'10.91"
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("D28").Value = "'9.08"
$ws.Range("E28").Value = '  -0.78%  '
$ws.Range("E29").Value = '  -2.74%  '
$ws.Range("E30").Value = '  -0.74%  '
$ws.Range("E31").Value = '  +0.69%  '
$ws.Range("E32").Value = '  +0.23%  '
$ws.Range("D33").Value = "'26.89"
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("E34").Value = '  +4.06%  '
$ws.Range("E35").Value = '  -1.74%  '
$ws.Range("D36").Value = '3.667.56'
$ws.Range("E36").Value = '  -0.23%  '
$ws.Range("D37").Value = "'8.45"
$ws.Range("E37").Value = '  +1.48%  '
$ws.Range("D39").Value = "'5.88"
$ws.Range("E39").Value = '  -5.95%  '
$ws.Range("D40").Value = "'178.91"
$ws.Range("E40").Value = '  +4.90%  '
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D43").Value = "'2.19"
$ws.Range("E43").Value = '  -1.64%  '
$ws.Range("E44").Value = '  -1.69%  '
$ws.Range("E45").Value = '  -1.73%  '
$ws.Range("E46").Value = '  +2.51%  '
$ws.Range("D47").Value = "'27.37"
$ws.Range("E47").Value = '  -3.00%  '
$ws.Range("E48").Value = '  -1.32%  '
$ws.Range("E49").Value = '  -3.51%  '
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("E51").Value = '  -3.87%  '
